# Refresh the live cryptocurrency price/volume snapshot (rows 2-51 of the
# "Price" and "Volume(1h)" columns) with the latest figures from the feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "58.995.24"
    "E2" = "  -0.39%  "
    "D3" = "2.511.48"
    "E3" = "  +0.49%  "
    "E4" = "  +0.07%  "
    "D5" = "534.52"
    "E5" = "  -0.32%  "
    "D6" = "135.77"
    "E6" = "  -1.29%  "
    "E7" = "  +0.30%  "
    "E8" = "  +0.21%  "
    "E9" = "  +0.24%  "
    "E10" = "  -1.10%  "
    "D11" = "5.43"
    "E11" = "  +1.26%  "
    "E12" = "  -0.35%  "
    "D13" = "2.956.36"
    "E13" = "  +0.35%  "
    "D14" = "58.899.25"
    "E14" = "  -0.03%  "
    "D15" = "22.83"
    "E16" = "  -1.11%  "
    "D17" = "2.502.51"
    "E17" = "  -0.27%  "
    "E18" = "  -0.28%  "
    "E19" = "  -0.33%  "
    "D20" = "322.93"
    "E20" = "  -0.73%  "
    "E22" = "  +0.90%  "
    "D23" = "65.11"
    "E23" = "  +0.64%  "
    "E24" = "  -0.12%  "
    "E25" = "  -0.97%  "
    "D26" = "0.998"
    "E26" = "  -0.92%  "
    "E27" = "  -0.73%  "
    "D28" = "0.0₃0765"
    "E28" = "  -1.34%  "
    "E29" = "  -3.03%  "
    "E30" = "  -1.39%  "
    "D31" = "169.78"
    "E31" = "  +1.53%  "
    "E33" = "  -4.05%  "
    "D34" = "1.37"
    "E34" = "  -2.51%  "
    "D35" = "18.41"
    "E35" = "  -0.90%  "
    "E36" = "  -1.93%  "
    "E37" = "  -3.12%  "
    "E38" = "  -1.60%  "
    "E39" = "  -4.03%  "
    "D40" = "282.36"
    "E40" = "  +0.45%  "
    "E41" = "  +0.47%  "
    "E42" = "  -5.03%  "
    "E43" = "  -0.02%  "
    "D44" = "129.53"
    "E44" = "  +1.24%  "
    "D45" = "10.90"
    "E45" = "  +0.21%  "
    "E46" = "  -0.49%  "
    "E47" = "  -2.07%  "
    "E48" = "  -2.37%  "
    "D49" = "17.31"
    "E49" = "  -0.51%  "
    "D50" = "1.761.02"
    "E50" = "  -0.61%  "
    "E51" = "  -0.49%  "
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $cell = $ws.Range($addr)

    # The Price column stores figures as literal text (so values like
    # "10.90" keep their trailing zero and big numbers like "58.995.24"
    # render using dots as thousand separators). Left alone, assigning a
    # plain numeric-looking string to a General-formatted cell makes Excel
    # silently reinterpret it as a number, so for those values we briefly
    # force Text formatting, write the string, then restore the original
    # style.
    if ($newValue -match '^[+-]?\d+(\.\d+)?$') {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = $origStyle
    } else {
        $cell.Value = $newValue
    }
}
